$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 372849
$ws.Range("E2").Value = 10719
$ws.Range("F2").Value = 10719
$ws.Range("G2").Value = 3841
$ws.Range("H2").Value = 4472
$ws.Range("I2").Value = 4472
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 467720
$ws.Range("L2").Value = 370476
$ws.Range("M2").Value = 97244
$ws.Range("N2").Value = 97244
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 4616
$ws.Range("Q2").Value = 20303
$ws.Range("R2").Value = -36754
$ws.Range("S2").Value = 15968
$ws.Range("T2").Value = 33878
$ws.Range("U2").Value = -13575
$ws.Range("V2").Value = 306431
$ws.Range("W2").Value = 2.88
$ws.Range("X2").Value = 1.2
$ws.Range("Y2").Value = 4.79
$ws.Range("Z2").Value = 0.99
$ws.Range("AA2").Value = 380.98
$ws.Range("AB2").Value = 1616.76
$ws.Range("AC2").Value = 4845
$ws.Range("AD2").Value = 10.23
$ws.Range("AE2").Value = 110962
$ws.Range("AF2").Value = 0.45
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 0.5
$ws.Range("AI2").Value = 4.9
$ws.Range("AJ2").Value = 92313000

# Row 3
$ws.Range("D3").Value = 260527
$ws.Range("E3").Value = 10078
$ws.Range("F3").Value = 10078
$ws.Range("G3").Value = 2768
$ws.Range("H3").Value = 3192
$ws.Range("I3").Value = 3192
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 423853
$ws.Range("L3").Value = 323284
$ws.Range("M3").Value = 100569
$ws.Range("N3").Value = 100569
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 4616
$ws.Range("Q3").Value = 54562
$ws.Range("R3").Value = -22466
$ws.Range("S3").Value = -32660
$ws.Range("T3").Value = 22724
$ws.Range("U3").Value = 31838
$ws.Range("V3").Value = 280527
$ws.Range("W3").Value = 3.87
$ws.Range("X3").Value = 1.23
$ws.Range("Y3").Value = 3.23
$ws.Range("Z3").Value = 0.72
$ws.Range("AA3").Value = 321.45
$ws.Range("AB3").Value = 1677.75
$ws.Range("AC3").Value = 3458
$ws.Range("AD3").Value = 10.66
$ws.Range("AE3").Value = 114757
$ws.Range("AF3").Value = 0.32
$ws.Range("AG3").Value = 170
$ws.Range("AH3").Value = 0.46
$ws.Range("AI3").Value = 4.67
$ws.Range("AJ3").Value = 92313000

# Row 4
$ws.Range("D4").Value = 211081
$ws.Range("E4").Value = 9982
$ws.Range("F4").Value = 9176
$ws.Range("G4").Value = -7690
$ws.Range("H4").Value = -6125
$ws.Range("I4").Value = -6130
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 400420
$ws.Range("L4").Value = 305690
$ws.Range("M4").Value = 94730
$ws.Range("N4").Value = 93603
$ws.Range("O4").Value = 1127
$ws.Range("P4").Value = 4616
$ws.Range("Q4").Value = 47693
$ws.Range("R4").Value = -20789
$ws.Range("S4").Value = -23766
$ws.Range("T4").Value = 18414
$ws.Range("U4").Value = 29279
$ws.Range("V4").Value = 259579
$ws.Range("W4").Value = 4.73
$ws.Range("X4").Value = -2.9
$ws.Range("Y4").Value = -6.31
$ws.Range("Z4").Value = -1.49
$ws.Range("AA4").Value = 322.69
$ws.Range("AB4").Value = 1546.1
$ws.Range("AC4").Value = -6641
$ws.Range("AD4").Value = -7.3
$ws.Range("AE4").Value = 106807
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 92313000

# Row 5
$ws.Range("D5").Value = 221722
$ws.Range("E5").Value = 10352
$ws.Range("F5").Value = 10352
$ws.Range("G5").Value = -12490
$ws.Range("H5").Value = -11917
$ws.Range("I5").Value = -12051
$ws.Range("J5").Value = 134
$ws.Range("K5").Value = 371394
$ws.Range("L5").Value = 289990
$ws.Range("M5").Value = 81404
$ws.Range("N5").Value = 78133
$ws.Range("O5").Value = 3271
$ws.Range("P5").Value = 4616
$ws.Range("Q5").Value = 25076
$ws.Range("R5").Value = -12259
$ws.Range("S5").Value = -12661
$ws.Range("T5").Value = 13516
$ws.Range("U5").Value = 11560
$ws.Range("V5").Value = 242495
$ws.Range("W5").Value = 4.67
$ws.Range("X5").Value = -5.38
$ws.Range("Y5").Value = -14.03
$ws.Range("Z5").Value = -3.09
$ws.Range("AA5").Value = 356.24
$ws.Range("AB5").Value = 1434.24
$ws.Range("AC5").Value = -13055
$ws.Range("AD5").Value = -3.26
$ws.Range("AE5").Value = 89155
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 92313000

# Row 6
$ws.Range("D6").Value = 261850
$ws.Range("E6").Value = 12769
$ws.Range("F6").Value = 12769
$ws.Range("G6").Value = 6490
$ws.Range("H6").Value = 5267
$ws.Range("I6").Value = 5094
$ws.Range("K6").Value = 396897
$ws.Range("L6").Value = 311917
$ws.Range("M6").Value = 84980
$ws.Range("N6").Value = 81928
$ws.Range("P6").Value = 4616
$ws.Range("Q6").Value = 571
$ws.Range("R6").Value = -12752
$ws.Range("S6").Value = 9942
$ws.Range("T6").Value = 11938
$ws.Range("U6").Value = -11367
$ws.Range("V6").Value = 261141
$ws.Range("W6").Value = 4.88
$ws.Range("X6").Value = 2.01
$ws.Range("Y6").Value = 6.36
$ws.Range("Z6").Value = 1.37
$ws.Range("AA6").Value = 367.05
$ws.Range("AB6").Value = 1628.52
$ws.Range("AC6").Value = 5518
$ws.Range("AD6").Value = 8.74
$ws.Range("AE6").Value = 93486
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 1360
$ws.Range("AH6").Value = 2.82
$ws.Range("AI6").Value = 23.4
$ws.Range("AJ6").Value = 92313000

# Row 7
$ws.Range("D7").Value = 250341
$ws.Range("E7").Value = 14002
$ws.Range("G7").Value = 7498
$ws.Range("H7").Value = 5332
$ws.Range("I7").Value = 5168
$ws.Range("K7").Value = 393426
$ws.Range("L7").Value = 306075
$ws.Range("M7").Value = 87351
$ws.Range("N7").Value = 84676
$ws.Range("P7").Value = 4619
$ws.Range("Q7").Value = 30098
$ws.Range("R7").Value = -13777
$ws.Range("S7").Value = -7205
$ws.Range("T7").Value = 13905
$ws.Range("U7").Value = 18352
$ws.Range("W7").Value = 5.59
$ws.Range("X7").Value = 2.13
$ws.Range("Y7").Value = 6.2
$ws.Range("Z7").Value = 1.35
$ws.Range("AA7").Value = 350.4
$ws.Range("AC7").Value = 5598
$ws.Range("AD7").Value = 5.81
$ws.Range("AE7").Value = 96621
$ws.Range("AF7").Value = 0.34
$ws.Range("AG7").Value = 1436
$ws.Range("AH7").Value = 4.42
$ws.Range("AI7").Value = 25.66

# Row 8
$ws.Range("D8").Value = 248734
$ws.Range("E8").Value = 14630
$ws.Range("G8").Value = 7868
$ws.Range("H8").Value = 5783
$ws.Range("I8").Value = 5704
$ws.Range("K8").Value = 392393
$ws.Range("L8").Value = 300803
$ws.Range("M8").Value = 91589
$ws.Range("N8").Value = 88758
$ws.Range("P8").Value = 4619
$ws.Range("Q8").Value = 24465
$ws.Range("R8").Value = -12616
$ws.Range("S8").Value = -3639
$ws.Range("T8").Value = 12385
$ws.Range("U8").Value = 14215
$ws.Range("W8").Value = 5.88
$ws.Range("X8").Value = 2.33
$ws.Range("Y8").Value = 6.58
$ws.Range("Z8").Value = 1.47
$ws.Range("AA8").Value = 328.43
$ws.Range("AC8").Value = 6179
$ws.Range("AD8").Value = 5.26
$ws.Range("AE8").Value = 101279
$ws.Range("AF8").Value = 0.32
$ws.Range("AG8").Value = 1552
$ws.Range("AH8").Value = 4.78
$ws.Range("AI8").Value = 25.12

# Row 9
$ws.Range("D9").Value = 252031
$ws.Range("E9").Value = 14743
$ws.Range("G9").Value = 8119
$ws.Range("H9").Value = 5974
$ws.Range("I9").Value = 5894
$ws.Range("K9").Value = 391819
$ws.Range("L9").Value = 295853
$ws.Range("M9").Value = 95967
$ws.Range("N9").Value = 92995
$ws.Range("P9").Value = 4619
$ws.Range("Q9").Value = 24440
$ws.Range("R9").Value = -13522
$ws.Range("S9").Value = -3937
$ws.Range("T9").Value = 13636
$ws.Range("U9").Value = 11758
$ws.Range("W9").Value = 5.85
$ws.Range("X9").Value = 2.37
$ws.Range("Y9").Value = 6.2
$ws.Range("Z9").Value = 1.52
$ws.Range("AA9").Value = 308.29
$ws.Range("AC9").Value = 6384
$ws.Range("AD9").Value = 5.09
$ws.Range("AE9").Value = 106114
$ws.Range("AF9").Value = 0.31
$ws.Range("AG9").Value = 1574
$ws.Range("AH9").Value = 4.84
$ws.Range("AI9").Value = 24.66
